$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
# "No" -> "ID"
$ws.Range("A1").Value = "ID"

# New "Kategori" header in E1, formatted like the other bold headers (copy D1's
# format which is bold Times New Roman) then set the font color explicitly to
# black (matches the new font added to the workbook: bold, Times New Roman,
# rgb 000000 instead of the theme color).
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("E1").Font.Color = 0
$ws.Range("E1").Value = "Kategori"

# --- New "Kategori" data column (E2:E23), formatted like column B ---
$ws.Range("B2:B23").Copy() | Out-Null
$ws.Range("E2:E23").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("E2").Value = "Hewani"
$ws.Range("E3").Value = "Hewani"
$ws.Range("E4").Value = "Nabati"
$ws.Range("E5").Value = "Hewani"
$ws.Range("E6").Value = "Nabati"
$ws.Range("E7").Value = "Nabati"
$ws.Range("E8").Value = "Nabati"
$ws.Range("E9").Value = "Hewani"
$ws.Range("E10").Value = "Hewani"
$ws.Range("E11").Value = "Hewani"
$ws.Range("E12").Value = "Hewan Air"
$ws.Range("E13").Value = "Nabati"
$ws.Range("E14").Value = "Nabati"
$ws.Range("E15").Value = "Hewani"
$ws.Range("E16").Value = "Hewan Air"
$ws.Range("E17").Value = "Hewan Air"
$ws.Range("E18").Value = "Hewan Air"
$ws.Range("E19").Value = "Hewan Air"
$ws.Range("E20").Value = "Hewan Air"
$ws.Range("E21").Value = "Hewan Air"
$ws.Range("E22").Value = "Hewan Air"
$ws.Range("E23").Value = "Hewan Air"

# --- Fix the last data row (Sarden Kaleng) ---
$ws.Range("A23").Value = "22"
$ws.Range("B23").Value = "Sarden Kaleng"

# --- Column E width ---
$ws.Range("E1").EntireColumn.ColumnWidth = 10.75

# --- Selection / active cell ---
$ws.Range("E7").Select() | Out-Null

Write-Output "done"
